# Regenerate the localization-status report: the "Status" column value used
# to read "Ready for handoff" wherever a file was still awaiting handoff;
# the refreshed report now shows "In Translation" instead. Update every
# sheet (Overview, zh-cn, de-de) and let the Status column narrow to fit
# the new (shorter) text, mirroring what Excel does when the column was
# sized to its contents.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# Overview sheet keeps the per-language status in columns E (zh-cn) and F (de-de).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1").ColumnWidth = 13.4101845877511
$overview.Range("F1").ColumnWidth = 13.4101845877511

# Each language sheet (zh-cn, de-de) keeps its own Status column in column C.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").ColumnWidth = 13.4101845877511

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").ColumnWidth = 13.4101845877511
